$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C12").Value = -11.59979999999999
$ws.Range("C27").Value = -12.6953
$ws.Range("C32").Value = -13.71020000000002
$ws.Range("C36").Value = -12.48340000000001
$ws.Range("C38").Value = -12.3425
$ws.Range("C46").Value = -14.79639999999999
$ws.Range("C54").Value = -12.90940000000001
$ws.Range("C55").Value = -13.33869999999999
$ws.Range("C56").Value = -12.42549999999999
$ws.Range("C67").Value = -11.0551
$ws.Range("C69").Value = -12.41279999999999
$ws.Range("C72").Value = -11.463
$ws.Range("C83").Value = -13.9567
$ws.Range("C86").Value = -13.23799999999999
$ws.Range("C91").Value = -10.4493
$ws.Range("C93").Value = -11.0337
$ws.Range("C99").Value = -12.964
